$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated numeric values for columns B:G across rows 2-7 (filtered save-game data)
$values = @{
    2 = @(0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0, 1.145820798638228)
    3 = @(0.1554434735375247, 0.05231270169004087, 3.082599426703578, 0.4998867070740569, 0, 3.790242309005201)
    4 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 0, 28.30127388105354)
    5 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
    6 = @(0.7287194209349384, 0.3375848360084654, 0.1529057820181812, 6.48142807727062, 0, 7.700638116232206)
    7 = @(0.1554434735375247, 0.05231270169004087, 3.082599426703578, 0.4998867070740569, 1, 3.790242309005201)
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = [char](66 + $i)  # B=66
        $ws.Range("$colLetter$row").Value = $cols[$i]
    }
}
